# "round 2 sea sky double" -- fill in the Round 2 ("N" column) bracket
# picks for the "Sea Beasties" category with the winners that advanced
# from Round 1, then leave the selection on the last cell touched (N32).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("N3").Value  = "Orca"
$ws.Range("N7").Value  = "Common Map Turtle"
$ws.Range("N11").Value = "Steller's Sea Eagle"
$ws.Range("N15").Value = "Hawaiian Monk Seal"
$ws.Range("N20").Value = "Walrus"
$ws.Range("N24").Value = "Pangolin"
$ws.Range("N28").Value = "Therapsid"
$ws.Range("N32").Value = "Swordfish"

$ws.Range("N32").Select()
